{"js": "// Add a new \"Seguir Programa\" row at the end of the tracking table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Append a new row; every cell except the first two is a simple \"-\" placeholder.\ntable.addRows(\"End\", 1, [[\n  \"Seguir Programa\",\n  \"\",\n  \"-\", \"-\", \"-\", \"-\", \"-\", \"-\", \"-\", \"-\", \"-\"\n]]);\nawait context.sync();\n\n// Make the \"Nombre CU\" cell of the new row bold (matches the style of the\n// other rows' first column).\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst newRowIndex = table.rowCount - 1;\nconst nameCell = table.getCell(newRowIndex, 0);\nconst nameRange = nameCell.getRange();\nnameRange.font.set({ bold: true });\nawait context.sync();\n\n// The \"Prueba 01 - Inicial\" cell's text was typed/pasted in three separate\n// pieces, so reproduce it as three distinct runs inside the paragraph\n// (rather than one merged run) by replacing the cell content with an OOXML\n// fragment that defines the runs explicitly.\nconst statusCell = table.getCell(newRowIndex, 1);\nconst statusRange = statusCell.getRange();\n\nconst statusOoxml =\n  '<?xml version=\"1.0\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n  '<w:r><w:t>No Realizada (listo para probar</w:t></w:r>' +\n  '<w:r><w:t>, Francisco o Fabricio</w:t></w:r>' +\n  '<w:r><w:t>)</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nstatusRange.insertOoxml(statusOoxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Add a new \"Seguir Programa\" row at the end of the tracking table.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$newRow = $tbl.Rows.Add()\n$rowIndex = $newRow.Index\n\n# Column 1 - \"Nombre CU\": bold label.\n$cell1 = $tbl.Cell($rowIndex, 1)\n$cell1.Range.Text = \"Seguir Programa\"\n$cell1.Range.Bold = 1\n\n# Column 2 - \"Prueba 01 - Inicial\" status. The text was authored as three\n# separate runs (same formatting, just typed/pasted separately), so insert\n# it as explicit OOXML runs instead of a single plain-text assignment.\n$cell2 = $tbl.Cell($rowIndex, 2)\n$xml2 = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr><w:r><w:t>No Realizada (listo para probar</w:t></w:r><w:r><w:t>, Francisco o Fabricio</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$cell2.Range.InsertXML($xml2)\n\n# Columns 3-11: placeholder \"-\" (no data yet for this new use case).\nfor ($c = 3; $c -le 11; $c++) {\n  $cell = $tbl.Cell($rowIndex, $c)\n  $cell.Range.Text = \"-\"\n}\n"}
